$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 5750.7407
$ws.Cells.Item(76, 9).Value = 4363.5
$ws.Cells.Item(76, 10).Value = 9714.286
$ws.Cells.Item(76, 11).Value = 4363.5
$ws.Cells.Item(76, 12).Value = 9714.286
$ws.Cells.Item(76, 13).Value = -4048.5
$ws.Cells.Item(76, 14).Value = -10344.286

$ws.Cells.Item(79, 8).Value = 5750.7407
$ws.Cells.Item(79, 9).Value = 4363.5
$ws.Cells.Item(79, 10).Value = 9714.286
$ws.Cells.Item(79, 11).Value = 4363.5
$ws.Cells.Item(79, 12).Value = 9714.286
$ws.Cells.Item(79, 13).Value = -3271.5
$ws.Cells.Item(79, 14).Value = -11898.286

$ws.Cells.Item(92, 8).Value = 4500.6665
$ws.Cells.Item(92, 9).Value = 5000.8
$ws.Cells.Item(92, 10).Value = 2000
$ws.Cells.Item(92, 11).Value = 5000.8
$ws.Cells.Item(92, 12).Value = 2000
$ws.Cells.Item(92, 13).Value = -3752.8
$ws.Cells.Item(92, 14).Value = -4496

$ws.Cells.Item(93, 8).Value = 34800.855
$ws.Cells.Item(93, 10).Value = 34800.855
$ws.Cells.Item(93, 12).Value = 34800.855
$ws.Cells.Item(93, 14).Value = -39792.855

$ws.Cells.Item(125, 8).Value = 1927.8572
$ws.Cells.Item(125, 9).Value = 2848.75
$ws.Cells.Item(125, 10).Value = 700
$ws.Cells.Item(125, 11).Value = 25638.75
$ws.Cells.Item(125, 12).Value = 6300
$ws.Cells.Item(125, 13).Value = -23178.75
$ws.Cells.Item(125, 14).Value = -11220

$ws.Cells.Item(126, 8).Value = 39000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 39000
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 39000
$ws.Cells.Item(126, 14).Value = -48880

$ws.Cells.Item(127, 8).Value = 934.37933
$ws.Cells.Item(127, 9).Value = 413.85715
$ws.Cells.Item(127, 10).Value = 1100
$ws.Cells.Item(127, 11).Value = 1241.57145
$ws.Cells.Item(127, 12).Value = 3300
$ws.Cells.Item(127, 13).Value = 3718.42855
$ws.Cells.Item(127, 14).Value = -13220

$ws.Cells.Item(128, 8).Value = 34867.5
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 34867.5
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 34867.5
$ws.Cells.Item(128, 14).Value = -44827.5

$ws.Cells.Item(129, 8).Value = 1250.9584
$ws.Cells.Item(129, 9).Value = 352.46155
$ws.Cells.Item(129, 10).Value = 1391.6868
$ws.Cells.Item(129, 11).Value = 1057.38465
$ws.Cells.Item(129, 12).Value = 4175.0604
$ws.Cells.Item(129, 13).Value = 3942.61535
$ws.Cells.Item(129, 14).Value = -14175.0604

$ws.Cells.Item(130, 8).Value = 36995
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 36995
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 36995
$ws.Cells.Item(130, 14).Value = -47035

$ws.Cells.Item(131, 8).Value = 2053.611
$ws.Cells.Item(131, 9).Value = 547
$ws.Cells.Item(131, 10).Value = 3936.875
$ws.Cells.Item(131, 11).Value = 1641
$ws.Cells.Item(131, 12).Value = 11810.625
$ws.Cells.Item(131, 13).Value = 3399
$ws.Cells.Item(131, 14).Value = -21890.625

$ws.Cells.Item(132, 8).Value = 5733.8647
$ws.Cells.Item(132, 9).Value = 3244.8928
$ws.Cells.Item(132, 10).Value = 13477.333
$ws.Cells.Item(132, 11).Value = 9734.678400000001
$ws.Cells.Item(132, 12).Value = 40431.999
$ws.Cells.Item(132, 13).Value = -7204.678400000001
$ws.Cells.Item(132, 14).Value = -45491.999

$ws.Cells.Item(133, 8).Value = 47589.875
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 47589.875
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 47589.875
$ws.Cells.Item(133, 14).Value = -57709.875

$ws.Cells.Item(134, 8).Value = 54701.8
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 54701.8
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 54701.8
$ws.Cells.Item(134, 14).Value = -64841.8

$ws.Cells.Item(135, 8).Value = 26317940
$ws.Cells.Item(135, 9).Value = 2341.7646
$ws.Cells.Item(135, 10).Value = 250000510
$ws.Cells.Item(135, 11).Value = 21075.8814
$ws.Cells.Item(135, 12).Value = 2250004590
$ws.Cells.Item(135, 13).Value = -18540.8814
$ws.Cells.Item(135, 14).Value = -2250009660

$ws.Cells.Item(136, 8).Value = 46421.285
$ws.Cells.Item(136, 9).Value = 10000
$ws.Cells.Item(136, 10).Value = 52491.5
$ws.Cells.Item(136, 11).Value = 10000
$ws.Cells.Item(136, 12).Value = 52491.5
$ws.Cells.Item(136, 13).Value = -4900
$ws.Cells.Item(136, 14).Value = -62691.5

$ws.Cells.Item(137, 8).Value = 3629.8076
$ws.Cells.Item(137, 9).Value = 1948.5555
$ws.Cells.Item(137, 10).Value = 7412.625
$ws.Cells.Item(137, 11).Value = 5845.666499999999
$ws.Cells.Item(137, 12).Value = 22237.875
$ws.Cells.Item(137, 13).Value = -3295.666499999999
$ws.Cells.Item(137, 14).Value = -27337.875

$ws.Cells.Item(138, 8).Value = 5210749
$ws.Cells.Item(138, 9).Value = 1257.8529
$ws.Cells.Item(138, 10).Value = 17862372
$ws.Cells.Item(138, 11).Value = 3773.5587
$ws.Cells.Item(138, 12).Value = 53587116
$ws.Cells.Item(138, 13).Value = 1366.4413
$ws.Cells.Item(138, 14).Value = -53597396

$ws.Cells.Item(139, 8).Value = 103500
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 103500
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 103500
$ws.Cells.Item(139, 14).Value = -113780

$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0

$ws.Cells.Item(141, 8).Value = 2321.077
$ws.Cells.Item(141, 9).Value = 2139.5
$ws.Cells.Item(141, 10).Value = 4500
$ws.Cells.Item(141, 11).Value = 6418.5
$ws.Cells.Item(141, 12).Value = 13500
$ws.Cells.Item(141, 13).Value = -1238.5
$ws.Cells.Item(141, 14).Value = -23860

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6611.0957
$ws.Cells.Item(32, 9).Value = 7235.731
$ws.Cells.Item(32, 10).Value = 5064.381
$ws.Cells.Item(32, 11).Value = 7235.731
$ws.Cells.Item(32, 12).Value = 5064.381
$ws.Cells.Item(32, 13).Value = -6948.731
$ws.Cells.Item(32, 14).Value = -5638.381

$ws.Cells.Item(102, 8).Value = 4880.45
$ws.Cells.Item(102, 9).Value = 5829.2144
$ws.Cells.Item(102, 10).Value = 2666.6667
$ws.Cells.Item(102, 11).Value = 5829.2144
$ws.Cells.Item(102, 12).Value = 2666.6667
$ws.Cells.Item(102, 13).Value = -4207.2144
$ws.Cells.Item(102, 14).Value = -5910.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 49250
$ws.Cells.Item(35, 10).Value = 49250
$ws.Cells.Item(35, 12).Value = 49250
$ws.Cells.Item(35, 14).Value = -49870

$ws.Cells.Item(105, 8).Value = 3801
$ws.Cells.Item(105, 9).Value = 1355
$ws.Cells.Item(105, 10).Value = 4412.5
$ws.Cells.Item(105, 11).Value = 1355
$ws.Cells.Item(105, 12).Value = 4412.5
$ws.Cells.Item(105, 13).Value = 392
$ws.Cells.Item(105, 14).Value = -7906.5

$ws.Cells.Item(134, 8).Value = 2706.5574
$ws.Cells.Item(134, 9).Value = 1821.72
$ws.Cells.Item(134, 10).Value = 6728.5454
$ws.Cells.Item(134, 11).Value = 5465.16
$ws.Cells.Item(134, 12).Value = 20185.6362
$ws.Cells.Item(134, 13).Value = -2930.16
$ws.Cells.Item(134, 14).Value = -25255.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1603.2222
$ws.Cells.Item(105, 9).Value = 857.25
$ws.Cells.Item(105, 10).Value = 2200
$ws.Cells.Item(105, 11).Value = 857.25
$ws.Cells.Item(105, 12).Value = 2200
$ws.Cells.Item(105, 13).Value = 889.75
$ws.Cells.Item(105, 14).Value = -5694

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2819
$ws.Cells.Item(34, 9).Value = 276
$ws.Cells.Item(34, 10).Value = 3666.6667
$ws.Cells.Item(34, 11).Value = 828
$ws.Cells.Item(34, 12).Value = 11000.0001
$ws.Cells.Item(34, 13).Value = -744
$ws.Cells.Item(34, 14).Value = -11168.0001

$ws.Cells.Item(39, 8).Value = 520
$ws.Cells.Item(39, 10).Value = 520
$ws.Cells.Item(39, 12).Value = 1560
$ws.Cells.Item(39, 14).Value = -2148

$ws.Cells.Item(55, 8).Value = 552.5
$ws.Cells.Item(55, 9).Value = 100
$ws.Cells.Item(55, 10).Value = 576.3158
$ws.Cells.Item(55, 11).Value = 300
$ws.Cells.Item(55, 12).Value = 1728.9474
$ws.Cells.Item(55, 13).Value = -123
$ws.Cells.Item(55, 14).Value = -2082.9474

$ws.Cells.Item(113, 8).Value = 902.325
$ws.Cells.Item(113, 9).Value = 794.56525
$ws.Cells.Item(113, 10).Value = 1048.1177
$ws.Cells.Item(113, 11).Value = 2383.69575
$ws.Cells.Item(113, 12).Value = 3144.3531
$ws.Cells.Item(113, 13).Value = -213.6957499999999
$ws.Cells.Item(113, 14).Value = -7484.3531

$ws.Cells.Item(122, 8).Value = 3423.75
$ws.Cells.Item(122, 9).Value = 10004
$ws.Cells.Item(122, 10).Value = 2483.7144
$ws.Cells.Item(122, 11).Value = 90036
$ws.Cells.Item(122, 12).Value = 22353.4296
$ws.Cells.Item(122, 13).Value = -87586
$ws.Cells.Item(122, 14).Value = -27253.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1407.8334
$ws.Cells.Item(97, 9).Value = 873.44446
$ws.Cells.Item(97, 10).Value = 3011
$ws.Cells.Item(97, 11).Value = 873.44446
$ws.Cells.Item(97, 12).Value = 3011
$ws.Cells.Item(97, 13).Value = -377.44446
$ws.Cells.Item(97, 14).Value = -4003

$ws.Cells.Item(132, 8).Value = 3554.6
$ws.Cells.Item(132, 9).Value = 2723.111
$ws.Cells.Item(132, 10).Value = 4530.696
$ws.Cells.Item(132, 11).Value = 8169.333
$ws.Cells.Item(132, 12).Value = 13592.088
$ws.Cells.Item(132, 13).Value = -5639.333
$ws.Cells.Item(132, 14).Value = -18652.088

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5415.1665
$ws.Cells.Item(40, 9).Value = 8688
$ws.Cells.Item(40, 10).Value = 4225.0454
$ws.Cells.Item(40, 11).Value = 8688
$ws.Cells.Item(40, 12).Value = 4225.0454
$ws.Cells.Item(40, 13).Value = -8552
$ws.Cells.Item(40, 14).Value = -4497.0454

$ws.Cells.Item(122, 8).Value = 5948.5
$ws.Cells.Item(122, 9).Value = 6900.25
$ws.Cells.Item(122, 10).Value = 5234.6875
$ws.Cells.Item(122, 11).Value = 20700.75
$ws.Cells.Item(122, 12).Value = 15704.0625
$ws.Cells.Item(122, 13).Value = -18250.75
$ws.Cells.Item(122, 14).Value = -20604.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 29954.857
$ws.Cells.Item(64, 10).Value = 29954.857
$ws.Cells.Item(64, 12).Value = 29954.857
$ws.Cells.Item(64, 14).Value = -30450.857

$ws.Cells.Item(67, 8).Value = 29954.857
$ws.Cells.Item(67, 10).Value = 29954.857
$ws.Cells.Item(67, 12).Value = 29954.857
$ws.Cells.Item(67, 14).Value = -31670.857
